# Add new columns I (I0) and J (IF) to Sheet1, populate header + data rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: copy formatting from existing header (H1) and set labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Data values for rows 2..61 (column I = "I0", column J = "IF")
$iValues = @(7,8,8,7,9,7,8,9,9,8,6,6,7,8,7,6,6,9,5,6,5,8,6,9,7,9,7,8,8,8,6,8,8,8,7,7,7,7,8,8,8,7,9,7,8,8,8,7,8,8,8,6,1,7,5,7,7,6,3,3)
$jValues = @(7,8,8,8,10,8,8,9,9,8,6,6,7,9,7,6,6,9,5,6,5,8,6,9,7,9,7,8,8,8,6,8,8,8,7,7,7,8,8,8,8,7,9,7,8,8,8,8,8,8,8,6,1,7,5,7,7,6,3,3)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}
